# Generate Report for Handoff
# Updates the localization-status workbook after a new handoff was generated
# for the acf528e8 and cf628c60 files: priority flips to machine-translation
# ("mt") and the handoff timestamps advance.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# New "Latest HO Xliff Generate Date" for the two re-handed-off files
$overview.Range("G3").Value = "2017-03-02 08:58:05"
$overview.Range("G4").Value = "2017-03-02 08:58:05"

# zh-cn sheet: Priority moved from "ht" to "mt"; Latest Handoff Datetime advanced
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2017-03-02 08:57:50"
$zhcn.Range("E4").Value = "mt"
$zhcn.Range("H4").Value = "2017-03-02 08:57:50"

# de-de sheet: Priority moved from "ht" to "mt"; Latest Handoff Datetime advanced
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2017-03-02 08:58:05"
$dede.Range("E4").Value = "mt"
$dede.Range("H4").Value = "2017-03-02 08:58:05"
